$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRINCIPAL")

$ws.Range("B9").Value = "DIO1659"
$ws.Range("H9").Value = "T - (Y 12/11/25_12H) - DF"

$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "12/11/25"
$ws.Range("I9").ClearFormats()

$ws.Range("J9").Value = "12H"
